$wb = $excel.ActiveWorkbook

# --- Sheet "2o Parcial", row 3 (student 6APV / TEMAS DE FILOSOFIA) ---
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("I3").Value = 6.6
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0

# --- Sheet "3er Parcial", row 3 (student 6APV / TEMAS DE FILOSOFIA) ---
$ws3 = $wb.Worksheets.Item("3er Parcial")
$ws3.Range("E3").Value = 20
$ws3.Range("F3").Value = 1
$ws3.Range("G3").Value = 95.23999999999999
$ws3.Range("H3").Value = 4.76
$ws3.Range("I3").Value = 7
